# Apply crypto price/volume updates as captured in the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.084.14"
$ws.Range("E2").Value = "  -2.62%  "
$ws.Range("D3").Value = "1.864.19"
$ws.Range("E3").Value = "  -2.40%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "306.12"
$ws.Range("E5").Value = "  -2.07%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.001"
$ws.Range("E6").Value = "  +0.15%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5151"
$ws.Range("E7").Value = "  -0.27%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3758"
$ws.Range("E8").Value = "  -0.53%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07156"
$ws.Range("E9").Value = "  -1.29%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8891"
$ws.Range("E10").Value = "  -1.81%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "20.68"
$ws.Range("E11").Value = "  -3.01%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07557"
$ws.Range("E12").Value = "  -1.31%  "
$ws.Range("D13").Value = "1.855.28"
$ws.Range("E13").Value = "  -3.16%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.307"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "89.63"
$ws.Range("E15").Value = "  -2.72%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.002"
$ws.Range("E16").Value = "  +0.02%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.000008473"
$ws.Range("E17").Value = "  -2.76%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "14.05"
$ws.Range("E18").Value = "  -3.26%  "
$ws.Range("E19").Value = "  +0.00%  "
$ws.Range("D20").Value = "27.128.98"
$ws.Range("E20").Value = "  -2.60%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.011"
$ws.Range("E21").Value = "  -2.76%  "
$ws.Range("D22").Value = "2.076.88"
$ws.Range("E22").Value = "  -3.56%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.442"
$ws.Range("E24").Value = "  -3.06%  "
$ws.Range("E25").Value = "  -1.71%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "145.62"
$ws.Range("E26").Value = "  -5.31%  "
$ws.Range("E27").Value = "  -2.20%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.087"
$ws.Range("E28").Value = "  -3.68%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "112.81"
$ws.Range("E29").Value = "  -1.89%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.660"
$ws.Range("E30").Value = "  -4.10%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.657"
$ws.Range("E31").Value = "  -4.25%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.09156"
$ws.Range("E32").Value = "  +0.76%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05100"
$ws.Range("E33").Value = "  -3.51%  "
$ws.Range("E34").Value = "  -3.48%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.154"
$ws.Range("E35").Value = "  -6.55%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7229"
$ws.Range("E36").Value = "  -7.26%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02036"
$ws.Range("E37").Value = "  -2.79%  "
$ws.Range("E38").Value = "  +0.20%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.485"
$ws.Range("E39").Value = "  -4.54%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.074"
$ws.Range("E40").Value = "  -1.79%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.5273"
$ws.Range("E41").Value = "  -5.75%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "6.465"
$ws.Range("E42").Value = "  -3.85%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "116.12"
$ws.Range("E43").Value = "  +0.52%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.272"
$ws.Range("E44").Value = "  -3.51%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.1463"
$ws.Range("E45").Value = "  -3.52%  "
$ws.Range("E46").Value = "  +0.22%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.4614"
$ws.Range("E47").Value = "  -4.43%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "9.943"
$ws.Range("E48").Value = "  -5.46%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.562"
$ws.Range("E49").Value = "  -3.46%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "36.51"
$ws.Range("E50").Value = "  -1.40%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "63.33"
$ws.Range("E51").Value = "  -5.49%  "
